$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pipe-delimited row data for the new bat 2289 sorting entries (rows 89-140).
# Columns: A=bat, B=date(serial), C=TT, D=thr, E=neg, F=lib_corr_thr,
#          G=min_win_sep, H=CD_thr, I=CD_nTT_thr, J(sorting comments)=K below
$rowData = @"
2289|43234|1|6|0|0.8|24|5|4|lots if firing! Thr by median is affected by the amount of spikes!
2289|43234|2|6|0|0.8|24|5|4|no cells
2289|43234|3|6|0|0.8|24|5|4|lots of activity but correlated across channels hard to isolate
2289|43234|4|6|0|0.8|24|5|4|lots of activity; hard to isolate
2289|43235|1|6|0|0.8|24|5|4|lots of activity; several interneurons
2289|43235|2|6|0|0.8|24|5|4|no cells
2289|43235|3|6|0|0.8|24|5|4|activity correlated across channels hard to isolate
2289|43235|4|6|0|0.8|24|5|4|activity correlated across channels hard to isolate
2289|43238|1|6|0|0.8|24|5|4|lots of activity (probably interneurons)
2289|43238|2|6|0|0.8|24|5|4|lots of clusters, some not very easy to isolate
2289|43238|3|6|1|0.8|24|5|4|no clear clusters (some inverted spikes hard to isolate)
2289|43238|4|6|0|0.8|24|5|4|no cells
2289|43240|1|6|0|0.8|24|5|4|only multi-units
2289|43240|2|6|0|0.8|24|5|4|no cells
2289|43240|3|7|0|0.8|24|5|4|many clusters
2289|43240|4|6|0|0.8|24|5|4|no cells
2289|43243|1|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43243|2|6|0|0.8|24|5|4|bad logger (+there are cells hard to cluster)
2289|43243|3|6|0|0.8|24|5|4|bad logger, but some MU can be sorted out
2289|43243|4|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43244|1|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43244|2|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43244|3|6|0|0.8|24|5|4|bad logger, but some MU can be sorted out
2289|43244|4|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43245|1|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43245|2|6|0|0.8|24|5|4|bad logger, but some MU can be sorted out
2289|43245|3|6|0|0.8|24|5|4|bad logger, but some MU can be sorted out
2289|43245|4|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43248|1|6|0|0.8|24|5|4|half of the recording with huge noise (+no cells)
2289|43248|2|6|0|0.8|24|5|4|half of the recording with huge noise
2289|43248|3|6|0|0.8|24|5|4|half of the recording with huge noise
2289|43248|4|6|0|0.8|24|5|4|half of the recording with huge noise (+no cells)
2289|43249|1|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43249|2|6|0|0.8|24|5|4|bad logger, but some MU can be sorted out
2289|43249|3|6|0|0.8|24|5|4|bad logger, some spikes but hard to isolate
2289|43249|4|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43251|1|6|0|0.8|24|5|4|bad logger, but some MU can be sorted out
2289|43251|2|6|0|0.8|24|5|4|bad logger, but some MU can be sorted out
2289|43251|3|6|0|0.8|24|5|4|bad logger, some spikes but hard to isolate
2289|43251|4|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43252|1|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43252|2|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43252|3|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43252|4|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43254|1|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43254|2|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43254|3|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43254|4|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43259|1|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43259|2|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43259|3|6|0|0.8|24|5|4|bad logger (+no cells)
2289|43259|4|6|0|0.8|24|5|4|bad logger (+no cells)
"@

$lines = $rowData -split "`n"

$startRow = 89
$rowIndex = 0
foreach ($line in $lines) {
    $f = $line -split "\|"
    $r = $startRow + $rowIndex

    $ws.Cells.Item($r, 1).Value = [int]$f[0]     # A bat
    $ws.Cells.Item($r, 2).Value = [int]$f[1]     # B date serial
    $ws.Cells.Item($r, 3).Value = [int]$f[2]     # C TT
    $ws.Cells.Item($r, 4).Value = [int]$f[3]     # D thr
    $ws.Cells.Item($r, 5).Value = [int]$f[4]     # E neg
    $ws.Cells.Item($r, 6).Value = [double]$f[5]  # F lib_corr_thr
    $ws.Cells.Item($r, 7).Value = [int]$f[6]     # G min_win_sep
    $ws.Cells.Item($r, 8).Value = [int]$f[7]     # H CD_thr
    $ws.Cells.Item($r, 9).Value = [int]$f[8]     # I CD_nTT_thr
    $ws.Cells.Item($r, 10).Value = "highpass"    # J filter_type
    $ws.Cells.Item($r, 11).Value = $f[9]         # K sorting comments

    $rowIndex = $rowIndex + 1
}

$lastRow = $startRow + $lines.Count - 1

# Match date formatting/style used by the existing date column (column B)
$ws.Range("B88").Copy()
$ws.Range("B" + $startRow + ":B" + $lastRow).PasteSpecial(-4122)

# Update view state to match: frozen pane scrolled down, new active selection
$ws.Application.CutCopyMode = $false
$ws.Range("A110").Select()
$excel.ActiveWindow.ScrollRow = 110
$ws.Range("I139").Select()
